$wb = $excel.ActiveWorkbook

# --- "Chart" sheet: append 8 new daily rows (2025-12-16 .. 2025-12-23) ---
$chart = $wb.Worksheets.Item("Chart")

$dates = @("2025-12-16","2025-12-17","2025-12-18","2025-12-19","2025-12-20","2025-12-21","2025-12-22","2025-12-23")
$notIndexed = @(200.0,200.0,200.0,200.0,200.0,200.0,200.0,200.0)
$indexed    = @(231.0,231.0,231.0,231.0,231.0,231.0,231.0,231.0)
$impressions = @(70.0,54.0,83.0,63.0,57.0,54.0,79.0,54.0)

$startRow = 74
for ($i = 0; $i -lt 8; $i++) {
    $r = $startRow + $i
    # Write the date as a formula first, then convert to a value in place via
    # copy/paste-special. A direct Value="2025-12-16" assignment gets smart-parsed
    # into a real date serial + date number format, which does not match the
    # source data (plain text dates, General format, same style as every other
    # cell in the sheet).
    $chart.Cells.Item($r, 1).Formula = '="' + $dates[$i] + '"'
    $chart.Cells.Item($r, 2).Value = $notIndexed[$i]
    $chart.Cells.Item($r, 3).Value = $indexed[$i]
    $chart.Cells.Item($r, 4).Value = $impressions[$i]
}
$dateRange = $chart.Range("A74:A81")
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)

# --- "Critical issues" sheet: refreshed reason/source/validation/pages rows ---
$critical = $wb.Worksheets.Item("Critical issues")

$critical.Cells.Item(2,1).Value = "Alternate page with proper canonical tag"
$critical.Cells.Item(2,2).Value = "Website"
$critical.Cells.Item(2,3).Value = "Failed"
$critical.Cells.Item(2,4).Value = 73.0

$critical.Cells.Item(3,1).Value = "Not found (404)"
$critical.Cells.Item(3,2).Value = "Website"
$critical.Cells.Item(3,3).Value = "Failed"
$critical.Cells.Item(3,4).Value = 56.0

$critical.Cells.Item(4,1).Value = "Page with redirect"
$critical.Cells.Item(4,2).Value = "Website"
$critical.Cells.Item(4,3).Value = "Failed"
$critical.Cells.Item(4,4).Value = 25.0

$critical.Cells.Item(5,1).Value = "Duplicate, Google chose different canonical than user"
$critical.Cells.Item(5,2).Value = "Google systems"
$critical.Cells.Item(5,3).Value = "Failed"
$critical.Cells.Item(5,4).Value = 19.0

$critical.Cells.Item(6,1).Value = "Crawled - currently not indexed"
$critical.Cells.Item(6,2).Value = "Google systems"
$critical.Cells.Item(6,3).Value = "Failed"
$critical.Cells.Item(6,4).Value = 8.0

$critical.Cells.Item(7,1).Value = "Excluded by ‘noindex’ tag"
$critical.Cells.Item(7,2).Value = "Website"
$critical.Cells.Item(7,3).Value = "Not Started"
$critical.Cells.Item(7,4).Value = 18.0

$critical.Cells.Item(8,1).Value = "Blocked by robots.txt"
$critical.Cells.Item(8,2).Value = "Website"
$critical.Cells.Item(8,3).Value = "Not Started"
$critical.Cells.Item(8,4).Value = 1.0

$critical.Cells.Item(9,1).Value = "Server error (5xx)"
$critical.Cells.Item(9,2).Value = "Website"
$critical.Cells.Item(9,3).Value = "Passed"
$critical.Cells.Item(9,4).Value = 0.0

$critical.Cells.Item(10,1).Value = "Discovered - currently not indexed"
$critical.Cells.Item(10,2).Value = "Google systems"
$critical.Cells.Item(10,3).Value = "Passed"
$critical.Cells.Item(10,4).Value = 0.0
